$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 5
$ws.Range("I22").Value = $true
$ws.Range("Z22").Value = "Dupertuis, Gaston"
$ws.Range("F23").Value = 11
$ws.Range("I23").Value = $false
$ws.Range("Z23").Value = "Nowotny, Martin"
$ws.Range("F27").Value = 24
$ws.Range("I27").Value = $true
$ws.Range("Z27").Value = "Aguirre, German"
$ws.Range("F28").Value = 73
$ws.Range("I28").Value = $false
$ws.Range("Z28").Value = "La Barba, Pablo"
$ws.Range("F37").Value = 15
$ws.Range("Z37").Value = "Rodriguez, Pablo"
$ws.Range("F38").Value = 31
$ws.Range("Z38").Value = "Fernandez, Carlos"
$ws.Range("F62").Value = 94
$ws.Range("Z62").Value = "Sartor, Yemel"
$ws.Range("F63").Value = 145
$ws.Range("Z63").Value = "Sueldo, Pablo"
$ws.Range("F80").Value = 158
$ws.Range("Z80").Value = "Campos, Dario"
$ws.Range("F81").Value = 160
$ws.Range("Z81").Value = "Chiara, Lucio"
$ws.Range("F92").Value = 16
$ws.Range("Z92").Value = "Rulfi, Daniel"
$ws.Range("F93").Value = 196
$ws.Range("Z93").Value = "Escobar, Esteban"
$ws.Range("F101").Value = 59
$ws.Range("Z101").Value = "Arrieta, Maximiliano"
$ws.Range("F102").Value = 124
$ws.Range("Z102").Value = "Mendieta, Elias"
$ws.Range("F103").Value = 127
$ws.Range("Z103").Value = "Badano, Pablo"
$ws.Range("F104").Value = 140
$ws.Range("Z104").Value = "Hagge, Pilar"
$ws.Range("F105").Value = 151
$ws.Range("Z105").Value = "Presel, Raul"
$ws.Range("F106").Value = 194
$ws.Range("Z106").Value = "Asenie, Santiago"
$ws.Range("F108").Value = 34
$ws.Range("I108").Value = $true
$ws.Range("K108").Value = 250
$ws.Range("O108").Value = 250
$ws.Range("S108").Value = 1
$ws.Range("W108").Value = "250 (S2022T01)"
$ws.Range("Z108").Value = "Levin, Raul"
$ws.Range("F109").Value = 87
$ws.Range("I109").Value = $false
$ws.Range("K109").Value = 0
$ws.Range("O109").Value = 0
$ws.Range("S109").Value = 0
$ws.Range("W109").Value = ""
$ws.Range("Z109").Value = "Perot, Martin"
$ws.Range("F126").Value = 223
$ws.Range("I126").Value = $false
$ws.Range("K126").Value = 0
$ws.Range("O126").Value = 0
$ws.Range("S126").Value = 0
$ws.Range("W126").Value = ""
$ws.Range("Z126").Value = "Aguer, Jose"
$ws.Range("F127").Value = 318
$ws.Range("I127").Value = $true
$ws.Range("K127").Value = 125
$ws.Range("O127").Value = 125
$ws.Range("S127").Value = 1
$ws.Range("W127").Value = "125 (S2022T01)"
$ws.Range("Z127").Value = "Bonelli, Marcos"
$ws.Range("F140").Value = 121
$ws.Range("Z140").Value = "Tenca, Javier"
$ws.Range("F141").Value = 211
$ws.Range("Z141").Value = "Miner, Alberto"
$ws.Range("F155").Value = 167
$ws.Range("Z155").Value = "Cossi, Francisco"
$ws.Range("F156").Value = 203
$ws.Range("Z156").Value = "Brian, Martin"
$ws.Range("F189").Value = 209
$ws.Range("Z189").Value = "Jose"
$ws.Range("F190").Value = 255
$ws.Range("Z190").Value = "Michea, Ignacio"
$ws.Range("F194").Value = 204
$ws.Range("Z194").Value = "Delgado, Pablo"
$ws.Range("F195").Value = 240
$ws.Range("Z195").Value = "Arrieta, Matias"
$ws.Range("F196").Value = 249
$ws.Range("Z196").Value = "Muller, Tomas"
$ws.Range("F197").Value = 256
$ws.Range("Z197").Value = "Portillo, Lucas"
$ws.Range("F220").Value = 292
$ws.Range("I220").Value = $true
$ws.Range("Z220").Value = "Dallinger, Humberto"
$ws.Range("F221").Value = 331
$ws.Range("I221").Value = $false
$ws.Range("Z221").Value = "Bracco, Fernando"
$ws.Range("F241").Value = 295
$ws.Range("Z241").Value = "Antunez, Pablo"
$ws.Range("F242").Value = 299
$ws.Range("Z242").Value = "Ferrero, Alejandro"
$ws.Range("F244").Value = 285
$ws.Range("I244").Value = $false
$ws.Range("Z244").Value = "Lell, Claudia"
$ws.Range("F245").Value = 304
$ws.Range("I245").Value = $true
$ws.Range("Z245").Value = "Velazquez, Pedro"
$ws.Range("F252").Value = 306
$ws.Range("Z252").Value = "Bertoli, Julian"
$ws.Range("F253").Value = 307
$ws.Range("Z253").Value = "Bertoli, Maximiliano"
